$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「早起きする方法」" (row 222) was removed from the spreadsheet.
# Deleting the entire row shifts every subsequent row up by one, which
# matches the renumbering seen across the rest of the sheet (old row 223
# becomes new row 222, ..., old row 270 becomes new row 269) and reduces
# the used range from A1:C270 to A1:C269.
$ws.Rows.Item(222).Delete()
